$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update TestCases value (B2): 40 -> 54
$ws.Range("B2").Value = "54"

# Update Instance value (D2): Automation2 -> Automation1
$ws.Range("D2").Value = "Automation1"

# Update the active selection to E2 (was C2)
$ws.Range("E2").Select()
